$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ZmsjB128"
$ws.Range("B2").Value = 23091535
$ws.Range("C2").Value = "ghqgnjf40"
$ws.Range("D2").Value = "d2t8TK#!"
$ws.Range("F2").Value = "JrmqPOzF"
$ws.Range("G2").Value = "aDfP"
